$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Controls")

# ------------------------------------------------------------------
# 1. Update the title - control count 58 -> 61
# ------------------------------------------------------------------
$ws.Range("A1").Value = "FSI Agent Governance Framework - All 61 Controls"

# ------------------------------------------------------------------
# 2. Insert a new row at position 49 for new Control 2.21
#    (shifts former rows 49-67 down to 50-68)
# ------------------------------------------------------------------
$ws.Rows("49:49").Insert()

# Data validation on C4:C56 auto-expanded to C4:C57 because of the
# inserted row; remove the validation from the extra trailing cell so
# the applied range shrinks back to the original C4:C56.
$ws.Range("C57").Validation.Delete()

# ------------------------------------------------------------------
# 3. Control 1.23 (row 27) - was a blank merged placeholder row;
#    fill in Name/Status and break the merge.
# ------------------------------------------------------------------
$ws.Range("A27:E27").UnMerge()
$ws.Range("B27").Value = "Step-Up Authentication for Agent Operations"
$ws.Range("C27").Value = "Not Started"

# ------------------------------------------------------------------
# 4. Control 2.19 (row 47) - was a blank merged placeholder row;
#    fill in Name/Status, clear the old placeholder formatting on
#    D/E and break the merge.
# ------------------------------------------------------------------
$ws.Range("A47:E47").UnMerge()
$ws.Range("B47:E47").ClearFormats()
$ws.Range("B47").Value = "Customer AI Disclosure and Transparency"
$ws.Range("C47").Value = "Not Started"
$ws.Range("D47:E47").ClearContents()

# ------------------------------------------------------------------
# 5. New Control 2.21 (row 49, created by the insert above)
# ------------------------------------------------------------------
$ws.Range("A49").NumberFormat = "@"
$ws.Range("A49").Value = "2.21"
$ws.Range("A49").ClearFormats()
$ws.Range("B49").Value = "AI Marketing Claims and Substantiation"
$ws.Range("C49").Value = "Not Started"

# ------------------------------------------------------------------
# 6. Control 3.9 (row 59 after the shift) - was a blank merged
#    placeholder row; fill in Name/Status and break the merge.
# ------------------------------------------------------------------
$ws.Range("A59:E59").UnMerge()
$ws.Range("B59").Value = "Agent Activity Dashboard"
$ws.Range("C59").Value = "Not Started"

# ------------------------------------------------------------------
# 7. Summary Dashboard sheet updates
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary Dashboard")
$ws2.Range("B4").Value = 23
$ws2.Range("B5").Value = 21
$ws2.Range("B8").Value = 61
